$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''52.335.40'
$ws.Range("E2").Value = '''  +1.43%  '

# Row 3
$ws.Range("D3").Value = '''2.889.55'
$ws.Range("E3").Value = '''  +3.86%  '

# Row 4
$ws.Range("E4").Value = '''  -0.04%  '

# Row 5
$ws.Range("D5").Value = '''352.55'
$ws.Range("E5").Value = '''  +0.13%  '

# Row 6
$ws.Range("D6").Value = '''113.23'
$ws.Range("E6").Value = '''  +3.94%  '

# Row 7
$ws.Range("D7").Value = '''0.561'
$ws.Range("E7").Value = '''  +1.94%  '

# Row 8
$ws.Range("E8").Value = '''  -0.03%  '

# Row 9
$ws.Range("D9").Value = '''0.625'
$ws.Range("E9").Value = '''  +2.72%  '

# Row 10
$ws.Range("D10").Value = '''40.64'
$ws.Range("E10").Value = '''  +2.53%  '

# Row 11
$ws.Range("D11").Value = '''0.136'
$ws.Range("E11").Value = '''  -0.35%  '

# Row 12
$ws.Range("D12").Value = '''0.0853'
$ws.Range("E12").Value = '''  +2.04%  '

# Row 13
$ws.Range("D13").Value = '''20.29'
$ws.Range("E13").Value = '''  +1.10%  '

# Row 14
$ws.Range("D14").Value = '''7.91'
$ws.Range("E14").Value = '''  +3.09%  '

# Row 15
$ws.Range("D15").Value = '''3.342.89'
$ws.Range("E15").Value = '''  +3.84%  '

# Row 16
$ws.Range("E16").Value = '''  +7.54%  '

# Row 17
$ws.Range("D17").Value = '''2.887.99'
$ws.Range("E17").Value = '''  +3.90%  '

# Row 18
$ws.Range("D18").Value = '''52.296.74'

# Row 19
$ws.Range("D19").Value = '''3.37'
$ws.Range("E19").Value = '''  +6.68%  '

# Row 20
$ws.Range("D20").Value = '''7.67'
$ws.Range("E20").Value = '''  -0.28%  '

# Row 21
$ws.Range("D21").Value = '''13.65'
$ws.Range("E21").Value = '''  +3.82%  '

# Row 22
$ws.Range("D22").Value = '''0.0₃0979'
$ws.Range("E22").Value = '''  +1.70%  '

# Row 23
$ws.Range("D23").Value = '''71.15'
$ws.Range("E23").Value = '''  +1.81%  '

# Row 24
$ws.Range("D24").Value = '''271.19'
$ws.Range("E24").Value = '''  +1.48%  '

# Row 25
$ws.Range("D25").Value = '''2.80'
$ws.Range("E25").Value = '''  +2.72%  '

# Row 26
$ws.Range("D26").Value = '''26.71'
$ws.Range("E26").Value = '''  +2.35%  '

# Row 27
$ws.Range("D27").Value = '''0.999'
$ws.Range("E27").Value = '''  +0.13%  '

# Row 28
$ws.Range("D28").Value = '''0.165'
$ws.Range("E28").Value = '''  +0.10%  '

# Row 29
$ws.Range("B29").Value = '''InjectiveProtocol'
$ws.Range("C29").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = '''39.10'
$ws.Range("E29").Value = '''  +5.75%  '

# Row 30
$ws.Range("B30").Value = '''Cosmos'
$ws.Range("C30").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '''10.60'
$ws.Range("E30").Value = '''  +3.58%  '

# Row 31
$ws.Range("D31").Value = '''2.28'
$ws.Range("E31").Value = '''  +1.79%  '

# Row 32
$ws.Range("E32").Value = '''  +2.60%  '

# Row 33
$ws.Range("D33").Value = '''52.76'
$ws.Range("E33").Value = '''  +2.08%  '

# Row 34
$ws.Range("D34").Value = '''5.84'
$ws.Range("E34").Value = '''  +2.60%  '

# Row 35
$ws.Range("D35").Value = '''0.0458'
$ws.Range("E35").Value = '''  +1.10%  '

# Row 36
$ws.Range("D36").Value = '''0.0900'
$ws.Range("E36").Value = '''  +8.09%  '

# Row 37
$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '''  -0.09%  '

# Row 38
$ws.Range("E38").Value = '''  +6.41%  '

# Row 39
$ws.Range("D39").Value = '''18.91'
$ws.Range("E39").Value = '''  +2.30%  '

# Row 40
$ws.Range("D40").Value = '''2.05'
$ws.Range("E40").Value = '''  +3.96%  '

# Row 41
$ws.Range("D41").Value = '''2.64'
$ws.Range("E41").Value = '''  +4.34%  '

# Row 42
$ws.Range("E42").Value = '''  +1.92%  '

# Row 43
$ws.Range("D43").Value = '''22.85'
$ws.Range("E43").Value = '''  +3.31%  '

# Row 44
$ws.Range("D44").Value = '''122.71'
$ws.Range("E44").Value = '''  +2.07%  '

# Row 45
$ws.Range("E45").Value = '''  +2.63%  '

# Row 46
$ws.Range("E46").Value = '''  +6.59%  '

# Row 47
$ws.Range("D47").Value = '''2.178.55'
$ws.Range("E47").Value = '''  +2.52%  '

# Row 48
$ws.Range("E48").Value = '''  +6.92%  '

# Row 49
$ws.Range("E49").Value = '''  +12.69%  '

# Row 50
$ws.Range("D50").Value = '''0.969'
$ws.Range("E50").Value = '''  +7.08%  '

# Row 51
$ws.Range("D51").Value = '''5.55'
$ws.Range("E51").Value = '''  +2.20%  '
